$d = $word.ActiveDocument

# Position a collapsed range immediately before the final paragraph mark
# (i.e. right after the last character of the document's last paragraph),
# then use InsertXML to splice in a brand-new list paragraph there.
$endPos = $d.Content.End
$r = $d.Range($endPos - 1, $endPos - 1)

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:pPr>' +
         '<w:pStyle w:val="ListParagraph"/>' +
         '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
       '</w:pPr>' +
       '<w:r><w:t xml:space="preserve">Please note that the therminol-66 temperature gets heated above 350 </w:t></w:r>' +
       '<w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>o</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve">C  (MAGNET_TEDS_HX_exit_Temp). This is due to the low flow rate in TEDS, as well as the simple </w:t></w:r>' +
       '<w:r><w:lastRenderedPageBreak/><w:t>model of the heat exchanger (using a constant UA). The heat exchanger model will need to be updated for more accurate representation of the system.</w:t></w:r>' +
       '</w:p>'

$r.InsertXML($xml) | Out-Null
